$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the existing row 82
# ("Fruta / hortaliza, semanal"). Insert a blank row so every
# subsequent record (old rows 82-150) shifts down by one (new rows
# 83-151), then populate the new row 82 with the latest reading.
$ws.Rows("82:82").Insert()

$ws.Cells.Item(82, 1).Value2 = 6
$ws.Cells.Item(82, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(82, 3).Value2 = "Metropolitana"
$ws.Cells.Item(82, 4).Value2 = 45280
$ws.Cells.Item(82, 5).Value2 = 13
$ws.Cells.Item(82, 6).Value2 = 100114007
$ws.Cells.Item(82, 7).Value2 = "Jengibre"
$ws.Cells.Item(82, 8).Value2 = "Sin especificar"
$ws.Cells.Item(82, 9).Value2 = "Primera"
$ws.Cells.Item(82, 10).Value2 = 100
$ws.Cells.Item(82, 11).Value2 = 25000
$ws.Cells.Item(82, 12).Value2 = 25000
$ws.Cells.Item(82, 13).Value2 = 25000
$ws.Cells.Item(82, 14).Value2 = "$/caja 13 kilos"
$ws.Cells.Item(82, 15).Value2 = "Perú"
$ws.Cells.Item(82, 16).Value2 = 1923
$ws.Cells.Item(82, 17).Value2 = 13
$ws.Cells.Item(82, 18).Value2 = "Hortaliza"
